$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current last row (row 60), pushing the
# "tooltip.SetActiveAtStartToggle" row down to row 63 along with its
# formatting (row height, styles, etc.)
$ws.Rows.Item(60).Insert()
$ws.Rows.Item(60).Insert()
$ws.Rows.Item(60).Insert()

# New translations for global (transform) attributes.
$ws.Range("A60").Value = "Position"
$ws.Range("B60").Value = "Position"
$ws.Range("C60").Value = "Posición"

$ws.Range("A61").Value = "Rotation"
$ws.Range("B61").Value = "Rotation"
$ws.Range("C61").Value = "Rotación"

$ws.Range("A62").Value = "Scale"
$ws.Range("B62").Value = "Scale"
$ws.Range("C62").Value = "Escala"

# Match the author's final selection / scroll position.
$ws.Range("B68").Select()
